# Weekly update: insert a new price row at the top of the data block (row 84),
# pushing all existing data rows down by one. The new row mirrors the layout
# of the other data rows for this market/product, with new Fecha/Volumen/
# Precio values for the latest week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 84 (and everything below it) down by one row, copying the
# formatting of row 84 into the newly inserted blank row.
$ws.Rows("84:84").Insert()

# Fill in the new row 84 with this week's data.
$ws.Range("A84").Value = 1
$ws.Range("B84").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C84").Value = "Arica y Parinacota"
$ws.Range("D84").Value = 45089
$ws.Range("E84").Value = 15
$ws.Range("F84").Value = 100112040
$ws.Range("G84").Value = "Cilantro"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 600
$ws.Range("K84").Value = 1800
$ws.Range("L84").Value = 2000
$ws.Range("M84").Value = 1917
$ws.Range("N84").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O84").Value = "Región de Arica y Parinacota"
$ws.Range("P84").Value = 958
$ws.Range("Q84").Value = 2
$ws.Range("R84").Value = "Hortaliza"
